$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dialog")

# Append a new dialog row (row 5): DialogID=4, Count=3, three sentences.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "您好，南京大学的招生老师！"
$ws.Range("D5").Value = "感谢您观看我的视频！"
$ws.Range("E5").Value = "这是让不懂编程的人做剧本用的。"

# Widen the sentence columns (C & D share one width, E and F get new widths too).
$ws.Range("C:D").ColumnWidth = 28
$ws.Range("E:E").ColumnWidth = 36.666666666667
$ws.Range("F:F").ColumnWidth = 11

# Move the live selection to match the author's final cursor position.
$ws.Range("E14").Select()
